# "Added Verification for Mingbo Inheritance document"
#
# The receiving-member verification columns (M: Verified, N: Status/Degree
# of Inheritance, O: Comments) are filled in for every inheritance-item row
# (rows 4-8) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$verified = "Yes, he explained somebit of it to me"
$degree   = "Know that where the code is  and how to run on my system, but still need to dive deep into its concept to get it fully understood."
$notes    = "Can explain it to the new members next semester and will try to understand maximum of it during summer vacation"

foreach ($r in 4..8) {
    $ws.Range("M$r").Value = $verified
    $ws.Range("N$r").Value = $degree
    $ws.Range("O$r").Value = $notes
}

# Row heights grew slightly to better fit the newly-wrapped text.
$ws.Rows.Item(4).RowHeight = 239.1
$ws.Rows.Item(5).RowHeight = 187.15
$ws.Rows.Item(6).RowHeight = 330
$ws.Rows.Item(8).RowHeight = 90

# Leave the cursor where the author finished typing.
$ws.Range("O4").Select() | Out-Null
